{"js": "// Helper: wrap a <w:body> inner-XML fragment into the minimal OOXML package\n// payload that Range.insertOoxml / Paragraph.insertOoxml expect.\nfunction wrapOoxml(bodyInner) {\n  return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyInner + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraphs we need by their current text content so the\n// script doesn't depend on brittle fixed indices.\nlet goBackParaIndex = -1;\nlet githubParaIndex = -1;\nlet commitParaIndex = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.trim() === \"\" && goBackParaIndex === -1 && i > 0) {\n    // candidate for the bookmark paragraph; confirm further below.\n  }\n  if (t.indexOf(\"GitHub is a code-hosting platform\") !== -1) {\n    githubParaIndex = i;\n  }\n  if (t.indexOf(\"Commit refers to\") !== -1) {\n    commitParaIndex = i;\n  }\n}\n// The bookmark paragraph is the empty paragraph immediately after the\n// title paragraph (\"Stephan Reyes\"), i.e. the second paragraph (index 1).\ngoBackParaIndex = 1;\n\n// 1) Strip the _GoBack bookmark from its original paragraph, leaving an\n//    empty paragraph behind (it gets re-added at the end of the document).\nconst goBackPara = paragraphs.items[goBackParaIndex];\ngoBackPara.insertOoxml(wrapOoxml(\"<w:p/>\"), Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Rewrite the \"3) GitHub is a code-hosting platform...\" paragraph,\n//    splitting it into more runs and wrapping the proper-noun words with\n//    <w:proofErr> spell-check markers (GitHub, BitKeeper, Bitbucket,\n//    SourceForge, Github's) without altering the visible text.\nconst githubParaOoxml =\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">3) </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>GitHub</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> is a code-hosting platform that allows for code to be saved online and shared with others for easy collaboration purposes. It was created in 2005 after the Linux kernel project had a falling out with the commercial company that created </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>BitKeeper</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>, seeing as they wanted to no longer market it as a free service. This caused the Linux development community, especially Linus Torvalds, to create a free to use code-hosting platform. This allows for anyone to publically create a repository for any/all other users to look at and modify, while keeping track of the changes created over time.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Other hosting services like </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Bitbucket</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> and </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>SourceForge</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> also exist, but neither have a community as large as that of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Github\\u2019s</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">. </w:t></w:r>' +\n  '</w:p>';\nparagraphs.items[githubParaIndex].insertOoxml(wrapOoxml(githubParaOoxml), Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Replace the unfinished \"Commit refers to \" list item with the\n//    completed item plus seven new Git-term list items (Push, A branch,\n//    Fork, Merge, Clone, Pull, Pull request).\nconst listItemsOoxml =\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Commit is used as saving </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">and displaying </w:t></w:r>' +\n    '<w:r><w:t>everything that has currently been added to the repository.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Push </w:t></w:r>' +\n    '<w:r><w:t>updates the remote repository with any commits made locally to a branch.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">A branch </w:t></w:r>' +\n    '<w:r><w:t>is used to reference all of the commits and allows for access to the history of the commits.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Fork refers to making a copy of the repository, usually used when collaborating with another person, in order to make changes but not change the core code.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Merge is used when one wants to combine two different branches to combine their changes. </w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Clone makes a copy of a project that exists remotely, however it is created locally.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Pull is used when one wants to update their local project with changes made in the remote project.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Pull request is used </w:t></w:r>' +\n    '<w:r><w:t>to tell those you are collaborating with about changes you have made to the remote repository.</w:t></w:r>' +\n  '</w:p>';\nparagraphs.items[commitParaIndex].insertOoxml(wrapOoxml(listItemsOoxml), Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Delete the (now two) trailing empty paragraphs that used to follow\n//    the \"Commit refers to \" item.\nconst tailParagraphs = body.paragraphs;\ntailParagraphs.load(\"items/text\");\nawait context.sync();\nconst trailingEmpty = [];\nfor (let i = tailParagraphs.items.length - 1; i >= 0; i--) {\n  if (tailParagraphs.items[i].text === \"\") {\n    trailingEmpty.push(tailParagraphs.items[i]);\n  } else {\n    break; // stop at the first non-empty paragraph walking backward\n  }\n}\n// Keep exactly one trailing empty paragraph to host the relocated\n// bookmark; delete any extras beyond that.\nfor (let i = 1; i < trailingEmpty.length; i++) {\n  trailingEmpty[i].delete();\n}\nawait context.sync();\n\n// 5) Re-create the _GoBack bookmark as the very last paragraph of the\n//    document (immediately before the section break).\nconst finalParagraphs = body.paragraphs;\nfinalParagraphs.load(\"items/text\");\nawait context.sync();\nconst lastPara = finalParagraphs.items[finalParagraphs.items.length - 1];\nconst bookmarkOoxml = '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>';\nlastPara.insertOoxml(wrapOoxml(bookmarkOoxml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n# Locate the paragraphs we need to touch by their current text content so\n# the script doesn't depend on brittle fixed indices.\n$goBackIndex = -1\n$githubIndex = -1\n$commitIndex = -1\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx++\n    $t = $p.Range.Text\n    if ($t -like \"*GitHub is a code-hosting platform*\") {\n        $githubIndex = $idx\n    }\n    if ($t -like \"*Commit refers to*\") {\n        $commitIndex = $idx\n    }\n}\n# The bookmark paragraph is the empty paragraph immediately after the\n# title paragraph (\"Stephan Reyes\"), i.e. the second paragraph.\n$goBackIndex = 2\n\n# 1) Strip the _GoBack bookmark from its original paragraph, leaving an\n#    empty paragraph behind (it gets re-added at the end of the document).\n$goBackPara = $d.Paragraphs($goBackIndex)\n$goBackPara.Range.InsertXML(\"<w:p $wNs/>\")\n\n# 2) Rewrite the \"3) GitHub is a code-hosting platform...\" paragraph,\n#    splitting it into more runs and wrapping the proper-noun words with\n#    <w:proofErr> spell-check markers (GitHub, BitKeeper, Bitbucket,\n#    SourceForge, Github's) without altering the visible text.\n$githubParaXml = \"<w:p $wNs>\" +\n    '<w:r><w:t xml:space=\"preserve\">3) </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>GitHub</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> is a code-hosting platform that allows for code to be saved online and shared with others for easy collaboration purposes. It was created in 2005 after the Linux kernel project had a falling out with the commercial company that created </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>BitKeeper</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>, seeing as they wanted to no longer market it as a free service. This caused the Linux development community, especially Linus Torvalds, to create a free to use code-hosting platform. This allows for anyone to publically create a repository for any/all other users to look at and modify, while keeping track of the changes created over time.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Other hosting services like </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Bitbucket</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> and </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>SourceForge</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> also exist, but neither have a community as large as that of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Github' + [char]0x2019 + 's</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">. </w:t></w:r>' +\n    '</w:p>'\n$githubPara = $d.Paragraphs($githubIndex)\n$githubPara.Range.InsertXML($githubParaXml)\n\n# 3) Replace the unfinished \"Commit refers to \" list item with the\n#    completed item plus seven new Git-term list items (Push, A branch,\n#    Fork, Merge, Clone, Pull, Pull request).\n$listParaOpen = \"<w:p $wNs><w:pPr><w:pStyle w:val=\"\"ListParagraph\"\"/><w:numPr><w:ilvl w:val=\"\"0\"\"/><w:numId w:val=\"\"2\"\"/></w:numPr></w:pPr>\"\n$listItemsXml =\n    $listParaOpen +\n        '<w:r><w:t xml:space=\"preserve\">Commit is used as saving </w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">and displaying </w:t></w:r>' +\n        '<w:r><w:t>everything that has currently been added to the repository.</w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t xml:space=\"preserve\">Push </w:t></w:r>' +\n        '<w:r><w:t>updates the remote repository with any commits made locally to a branch.</w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t xml:space=\"preserve\">A branch </w:t></w:r>' +\n        '<w:r><w:t>is used to reference all of the commits and allows for access to the history of the commits.</w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t>Fork refers to making a copy of the repository, usually used when collaborating with another person, in order to make changes but not change the core code.</w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t xml:space=\"preserve\">Merge is used when one wants to combine two different branches to combine their changes. </w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t>Clone makes a copy of a project that exists remotely, however it is created locally.</w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t>Pull is used when one wants to update their local project with changes made in the remote project.</w:t></w:r>' +\n    '</w:p>' +\n    $listParaOpen +\n        '<w:r><w:t xml:space=\"preserve\">Pull request is used </w:t></w:r>' +\n        '<w:r><w:t>to tell those you are collaborating with about changes you have made to the remote repository.</w:t></w:r>' +\n    '</w:p>'\n$commitPara = $d.Paragraphs($commitIndex)\n$commitPara.Range.InsertXML($listItemsXml)\n\n# 4) Delete the (now two) trailing empty paragraphs that used to follow\n#    the \"Commit refers to \" item, keeping exactly one to host the\n#    relocated bookmark.\n$trailing = @()\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq \"\") {\n        $trailing += $i\n    } else {\n        break\n    }\n}\nfor ($i = 0; $i -lt ($trailing.Count - 1); $i++) {\n    $d.Paragraphs($trailing[$i]).Range.Delete()\n}\n\n# 5) Re-create the _GoBack bookmark as the very last paragraph of the\n#    document (immediately before the section break).\n$lastIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($lastIndex)\n$bookmarkXml = \"<w:p $wNs><w:bookmarkStart w:id=\"\"0\"\" w:name=\"\"_GoBack\"\"/><w:bookmarkEnd w:id=\"\"0\"\"/></w:p>\"\n$lastPara.Range.InsertXML($bookmarkXml)\n\nWrite-Output \"done\"\n"}
